# Cash Invoice workbook update:
#  - add three new invoice sheets (Ambika Jewellers 004, Nighale 005, Adv Manerkar 006)
#  - move the active tab to the new "Nighale 005" sheet
#  - leave a new selection behind on "Shezwan House 003"

$wb = $excel.ActiveWorkbook

$xlCenter = -4108
$xlNone = -4142
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10

function Format-HeaderCell($cell) {
    $cell.Font.Bold = $true
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 10
    $cell.HorizontalAlignment = $xlCenter
    $cell.VerticalAlignment = $xlCenter
    $cell.WrapText = $true
    $cell.Borders.LineStyle = 1
}

function Format-BodyCell($cell) {
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 10
    $cell.HorizontalAlignment = $xlCenter
    $cell.VerticalAlignment = $xlCenter
    $cell.WrapText = $true
    $cell.Borders.LineStyle = 1
}

# Total-row cell whose box only needs a left edge (used for the left-most cell
# of a merged "TOTAL"/"Total" label run).
function Format-TotalLeftCell($cell) {
    $cell.Font.Bold = $true
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 10
    $cell.HorizontalAlignment = $xlCenter
    $cell.VerticalAlignment = $xlCenter
    $cell.WrapText = $true
    $cell.Borders.LineStyle = 1
    $cell.Borders.Item($xlEdgeRight).LineStyle = $xlNone
}

# Total-row cell with no left/right edges (interior of the merged label run).
function Format-TotalMidCell($cell) {
    $cell.Font.Bold = $true
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 10
    $cell.HorizontalAlignment = $xlCenter
    $cell.VerticalAlignment = $xlCenter
    $cell.WrapText = $true
    $cell.Borders.LineStyle = 1
    $cell.Borders.Item($xlEdgeLeft).LineStyle = $xlNone
    $cell.Borders.Item($xlEdgeRight).LineStyle = $xlNone
}

# Total-row cell (the grand-total amount) keeping the full box.
function Format-TotalRightCell($cell) {
    $cell.Font.Bold = $true
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 10
    $cell.HorizontalAlignment = $xlCenter
    $cell.VerticalAlignment = $xlCenter
    $cell.WrapText = $true
    $cell.Borders.LineStyle = 1
}

function Add-InvoiceSheet($name) {
    $after = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add([System.Type]::Missing, $after)
    $ws.Name = $name
    $ws.Tab.Color = 5296274   # FF92D050
    return $ws
}

# ---------------------------------------------------------------------------
# Sheet 4: Ambika Jewellers 004
# ---------------------------------------------------------------------------
$ws4 = Add-InvoiceSheet "Ambika Jewellers 004"

$ws4.Range("A1").Value = "SR NO"
$ws4.Range("B1").Value = "ITEM DESCRIPTION"
$ws4.Range("C1").Value = "QTY"
$ws4.Range("D1").Value = "PRICE"
$ws4.Range("E1").Value = "AMOUNT"
foreach ($col in @("A","B","C","D","E")) { Format-HeaderCell $ws4.Range($col + "1") }
$ws4.Rows.Item(1).RowHeight = 25.8

$ws4.Range("A2").Value = 1
$ws4.Range("B2").Value = "CCTV Cable 3+1"
$ws4.Range("C2").Value = 10
$ws4.Range("D2").Value = 90
$ws4.Range("E2").Formula = "=C2*D2"

$ws4.Range("A3").Value = 2
$ws4.Range("B3").Value = "Service Calls Charges"
$ws4.Range("C3").Value = 1
$ws4.Range("D3").Value = 900
$ws4.Range("E3").Formula = "=C3*D3"

foreach ($r in 2..3) {
    foreach ($col in @("A","B","C","D","E")) { Format-BodyCell $ws4.Range($col + $r) }
}

$ws4.Range("A4").Value = "TOTAL"
$ws4.Range("E4").Formula = "=SUM(E2:E3)"
Format-TotalLeftCell $ws4.Range("A4")
Format-TotalMidCell $ws4.Range("B4")
Format-TotalMidCell $ws4.Range("C4")
Format-TotalMidCell $ws4.Range("D4")
Format-TotalRightCell $ws4.Range("E4")
$ws4.Range("A4:D4").Merge()

$ws4.Columns.Item(2).ColumnWidth = 22.44

$ws4.Range("J24").Select()

# ---------------------------------------------------------------------------
# Sheet 5: Nighale 005
# ---------------------------------------------------------------------------
$ws5 = Add-InvoiceSheet "Nighale 005"

$ws5.Range("A1").Value = "SR NO"
$ws5.Range("B1").Value = "ITEM DESCRIPTION"
$ws5.Range("C1").Value = "QTY"
$ws5.Range("D1").Value = "PRICE"
$ws5.Range("E1").Value = "AMOUNT"
foreach ($col in @("A","B","C","D","E")) { Format-HeaderCell $ws5.Range($col + "1") }
$ws5.Rows.Item(1).RowHeight = 41.4

$ws5.Range("A2").Value = 1
$ws5.Range("B2").Value = "Service call charges"
$ws5.Range("C2").Value = 1
$ws5.Range("D2").Value = 1000
$ws5.Range("E2").Formula = "=C2*D2"
foreach ($col in @("A","B","C","D","E")) { Format-BodyCell $ws5.Range($col + "2") }
$ws5.Rows.Item(2).RowHeight = 27.6

$ws5.Range("A3").Value = "Total"
$ws5.Range("E3").Formula = "=SUM(E2)"
Format-TotalLeftCell $ws5.Range("A3")
Format-TotalMidCell $ws5.Range("B3")
Format-TotalMidCell $ws5.Range("C3")
Format-TotalMidCell $ws5.Range("D3")
Format-TotalRightCell $ws5.Range("E3")
$ws5.Range("A3:D3").Merge()

$ws5.Columns.Item(2).ColumnWidth = 14.78

# ---------------------------------------------------------------------------
# Sheet 6: Adv Manerkar 006
# ---------------------------------------------------------------------------
$ws6 = Add-InvoiceSheet "Adv Manerkar 006"

$ws6.Range("A1").Value = "SR NO"
$ws6.Range("B1").Value = "ITEM DESCRIPTION"
$ws6.Range("C1").Value = "QTY"
$ws6.Range("D1").Value = "PRICE"
$ws6.Range("E1").Value = "AMOUNT"
foreach ($col in @("A","B","C","D","E")) { Format-HeaderCell $ws6.Range($col + "1") }

$ws6.Range("A2").Value = 1
$ws6.Range("B2").Value = "Exide Battery 150 AH"
$ws6.Range("C2").Value = 1
$ws6.Range("D2").Value = 14000
$ws6.Range("E2").Formula = "=C2*D2"
foreach ($col in @("A","B","C","D","E")) { Format-BodyCell $ws6.Range($col + "2") }

$ws6.Range("A3").Value = "Total"
$ws6.Range("E3").Formula = "=SUM(E2)"
Format-TotalLeftCell $ws6.Range("A3")
Format-TotalMidCell $ws6.Range("B3")
Format-TotalMidCell $ws6.Range("C3")
Format-TotalMidCell $ws6.Range("D3")
Format-TotalRightCell $ws6.Range("E3")
$ws6.Range("A3:D3").Merge()

$ws6.Columns.Item(2).ColumnWidth = 26

$ws6.Range("A1:E3").Select()

# ---------------------------------------------------------------------------
# Selection / active-tab bookkeeping
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Shezwan House 003")
$ws3.Range("B8").Select()

$ws5.Activate()
$ws5.Range("B2").Select()
